$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1): update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2140
$ws1.Range("F5").Value = 11171
$ws1.Range("F9").Value = 214
$ws1.Range("F10").Value = 11079
$ws1.Range("F12").Value = 1142
$ws1.Range("F13").Value = 42
$ws1.Range("F15").Value = 5557
$ws1.Range("F17").Value = 3437

# Sheet "全部类型" (rId4 / sheet4): same underlying events, mirrored counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2140
$ws4.Range("F7").Value = 11171
$ws4.Range("F11").Value = 214
$ws4.Range("F12").Value = 11079
$ws4.Range("F14").Value = 1142
$ws4.Range("F15").Value = 42
$ws4.Range("F17").Value = 5557
$ws4.Range("F19").Value = 3437
